$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Step 1: insert rows to make room for new entries, working from the bottom up
# so that earlier insert operations do not shift the row numbers used by later ones.

# Insert a single blank row before old row 23 ("These are for OUTPUT..."),
# which will become the blank row 30 separating row 29 (athletetest) from row 31.
$ws.Rows.Item(23).EntireRow.Insert()

# Insert 2 rows before old row 22 (athletetest) for TopopprankIndiv / TopopprankRelay (new rows 27-28)
$ws.Rows.Item(22).Resize(2).EntireRow.Insert()

# Insert 2 rows before old row 17 (Maxevent) for indivplcscore / relayplcscore (new rows 20-21)
$ws.Rows.Item(17).Resize(2).EntireRow.Insert()

# Insert 3 rows before old row 11 (homerank) for playperf / playperfstart / playperfMR (new rows 11-13)
$ws.Rows.Item(11).Resize(3).EntireRow.Insert()

# Step 2: (re)write all cell contents for rows 1-29 and 31-35 to match the final layout.
# Clear D12 and D13 explicitly since those rows have no "Use" text.

$ws.Range("A1").Value = "Variable"
$ws.Range("B1").Value = "Type"
$ws.Range("C1").Value = "Index"
$ws.Range("D1").Value = "Use"

$ws.Range("A2").Value = "MaxSolveTime"
$ws.Range("B2").Value = "int"
$ws.Range("C2").Value = "N/A"
$ws.Range("D2").Value = "Max number of minutes the optimization will run, if optimality conditions are not met, before returning the best found feasible solution"

$ws.Range("A3").Value = "OptGap"
$ws.Range("B3").Value = "float"
$ws.Range("C3").Value = "N/A"
$ws.Range("D3").Value = "Optimality condition, the decimal (0-1) gap to satisfy the optimality conditions"

$ws.Range("A4").Value = "ActAthNum"
$ws.Range("B4").Value = "int"
$ws.Range("C4").Value = "N/A"
$ws.Range("D4").Value = "Number of actual athletes on the roster"

$ws.Range("A5").Value = "athleteFull"
$ws.Range("B5").Value = "list"
$ws.Range("C5").Value = "int"
$ws.Range("D5").Value = "list of Names of each of the athletes (assumed unique). This could be changed to swimmerIDs, but have a lookup table to connect to the actual name"

$ws.Range("A6").Value = "scenario"
$ws.Range("B6").Value = "list"
$ws.Range("C6").Value = "int"
$ws.Range("D6").Value = "Names of the different scenarios will be the index for other dictionaries"

$ws.Range("A7").Value = "indiv"
$ws.Range("B7").Value = "list"
$ws.Range("C7").Value = "int"
$ws.Range("D7").Value = "Names of all the individual events"

$ws.Range("A8").Value = "event_noMR"
$ws.Range("B8").Value = "list"
$ws.Range("C8").Value = "int"
$ws.Range("D8").Value = "Names of the events EXCLUDING the medley relay strokes."

$ws.Range("A9").Value = "relaynoMR"
$ws.Range("B9").Value = "list"
$ws.Range("C9").Value = "int"
$ws.Range("D9").Value = "Names of the relay events EXCLUDING the medley relay strokes"

$ws.Range("A10").Value = "stroke"
$ws.Range("B10").Value = "list"
$ws.Range("C10").Value = "int"
$ws.Range("D10").Value = "Names of strokes/legs in Medley relay"

$ws.Range("A11").Value = "playperf"
$ws.Range("B11").Value = "list"
$ws.Range("C11").Value = "athlete, indiv"
$ws.Range("D11").Value = "pred times for each athlete in each individual event"

$ws.Range("A12").Value = "playperfstart"
$ws.Range("B12").Value = "list"
$ws.Range("C12").Value = "athlete, relaynoMR"

$ws.Range("A13").Value = "playperfMR"
$ws.Range("B13").Value = "list"
$ws.Range("C13").Value = "athlete, stroke"

$ws.Range("A14").Value = "homerank"
$ws.Range("B14").Value = "list"
$ws.Range("C14").Value = "int"
$ws.Range("D14").Value = "Numbers 1,2,3"

$ws.Range("A15").Value = "event11"
$ws.Range("B15").Value = "list"
$ws.Range("C15").Value = "int"
$ws.Range("D15").Value = "The list of actual events contested: each individual event, each relay event (not legs - the actual events)"

$ws.Range("A16").Value = "place"
$ws.Range("B16").Value = "list"
$ws.Range("C16").Value = "int"
$ws.Range("D16").Value = "The possible finishing places for a race. Typically: 1,2,..,8"

$ws.Range("A17").Value = "event"
$ws.Range("B17").Value = "list"
$ws.Range("C17").Value = "int"
$ws.Range("D17").Value = "the list of individual events, relay events, MR strokes. ""200MR"" does double duty as the label for the whole MR and to identify being assigned to the first leg of the 200MR."

$ws.Range("A18").Value = "EventNoTimeArray"
$ws.Range("B18").Value = "list"
$ws.Range("C18").Value = "int"
$ws.Range("D18").Value = "List of max times used for athletes with no listed time in event (I don't think this is currently used - we assume the input data has already done this)"

$ws.Range("A19").Value = "scenprob"
$ws.Range("B19").Value = "list"
$ws.Range("C19").Value = "int"
$ws.Range("D19").Value = "List of probabilities for each scenario. In current version this will be the Nash equilibrium likelihoods passed from the game theory analysis for the opponent's possible lineups."

$ws.Range("A20").Value = "indivplcscore"
$ws.Range("B20").Value = "list"
$ws.Range("C20").Value = "place"
$ws.Range("D20").Value = "List of points for finishing places in INDIVIDUAL events"

$ws.Range("A21").Value = "relayplcscore"
$ws.Range("B21").Value = "list"
$ws.Range("C21").Value = "place"
$ws.Range("D21").Value = "List of points for finishing places in RELAY events"

$ws.Range("A22").Value = "Maxevent"
$ws.Range("B22").Value = "int"
$ws.Range("D22").Value = "Maximum number of TOTAL events per athlete"

$ws.Range("A23").Value = "Maxrelayevent"
$ws.Range("B23").Value = "int"
$ws.Range("D23").Value = "Maximum number of RELAY events an athlete can be assigned to"

$ws.Range("A24").Value = "MaxIndevent"
$ws.Range("B24").Value = "int"
$ws.Range("D24").Value = "Maximum number of INDIVIDUAL events an athlete can be assigned to"

$ws.Range("A25").Value = "BigM"
$ws.Range("B25").Value = "list int"
$ws.Range("C25").Value = "event11"
$ws.Range("D25").Value = "If there are no athletes assigned to an event (empty slots on relay, can't feasibly fill all three slots), we need to give the time ""rvar"" a large enough value that it optimization doesn't treat it like it ""won"" the event. There are computational costs for using values that are too large. BigM should only be as large as necessary!"

$ws.Range("A26").Value = "Tophomerank"
$ws.Range("B26").Value = "int"
$ws.Range("D26").Value = "Currently 3, just the number of possible ranks for a hometeam athlete (if only two were allowed, this would be 2)."

$ws.Range("A27").Value = "TopopprankIndiv"
$ws.Range("B27").Value = "int"
$ws.Range("D27").Value = "Currently 3, the number of opponents assigned to each INDIVIDUAL. In non-dual meets, you could have as many as 5 other opponents in a race"

$ws.Range("A28").Value = "TopopprankRelay"
$ws.Range("B28").Value = "int"
$ws.Range("D28").Value = "Currently 3, the number of opponents in each relay event. In non-dual meets, you could have as many as 5 other opponents in a race"

$ws.Range("A29").Value = "athletetest"
$ws.Range("B29").Value = "string"
$ws.Range("D29").Value = "This is for VORP and tells the model if you want to include all the athletes for the all the events. Only needed if we uncomment the VORP/WAR section"

$ws.Range("A31").Value = "These are for OUTPUT in the Excel sheet, but something similary (or identical) will need to be recorded."

$ws.Range("A32").Value = "HomeAthPredTime"
$ws.Range("B32").Value = "list string"
$ws.Range("C32").Value = "homerank, event11"
$ws.Range("D32").Value = "Will hold the predicted TIME of the 1st, 2nd, 3rd assigned home athlete/team in each event"

$ws.Range("A33").Value = "HomeAthFinPlace"
$ws.Range("B33").Value = "list int"
$ws.Range("C33").Value = "homerank, event11"
$ws.Range("D33").Value = "Will hold the predicted FINISH PLACE of the 1st, 2nd, 3rd assigned home athlete/team in each event"

$ws.Range("A34").Value = "HomeAthAssgnNamesIndiv"
$ws.Range("B34").Value = "list "
$ws.Range("C34").Value = "indiv, athlete"
$ws.Range("D34").Value = "Will hold the NAMES of the 1st, 2nd, 3rd assigned home athlete in each INDIVIDUAL event"

$ws.Range("A35").Value = "HomeAthAssgnNamesRelay"
$ws.Range("B35").Value = "list"
$ws.Range("C35").Value = "(0,1,2,3), (0,1,2,…,8)"
$ws.Range("D35").Value = "Creates an array of the names of the athletes doing each relay (9 relays, 4 athletes in each). Ordered by leadoff legs in non-MR and strokes for MR. Need to see the structure of the worksheet ""4. Assignment and Prediction"" for this to make sense."

# Sheet 2 change: "Print statements..." note index changed due to shared-string reordering,
# but its displayed text is identical.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B7").Value = "Print statements need to be updated for Python 3"

# Match the final selected cell recorded in the workbook view.
$ws.Range("D28").Select() | Out-Null
